$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UTIJoinville")

$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = 62
$ws.Cells.Item(46, 3).Value = 2
$ws.Cells.Item(46, 4).Value = 12
$ws.Cells.Item(46, 5).Value = 29
$ws.Cells.Item(46, 6).Value = 76
$ws.Cells.Item(46, 7).Value = 105
